# Updated arch diagram and bullets
#
# Re-positions / re-sizes several shapes on slide 1 of the Cloud One
# Conformity architecture diagram, and turns on the missing arrowhead for
# one of the connectors.
#
# Note: PowerPoint's COM object model stores Shape.Left/Top/Width/Height in
# points, while the underlying OOXML stores offsets/extents in EMU
# (1 pt = 12700 EMU). The point values below were chosen so that they land
# safely in the middle of the EMU bucket that rounds/truncates back to the
# exact target EMU value from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 2 - "Rectangle 123" (id 124, "AWS Cloud" outer box)
#   off  1794514,640081 -> 1794515,640081
#   ext  6365417,4536478 -> 6374126,4881154
$sh = $s.Shapes.Item(2)
$sh.Left   = 141.30043030078522
$sh.Top    = 50.400117840157925
$sh.Width  = 501.89971923937776
$sh.Height = 384.34288024567616

# Shape 7 - "TextBox 29" (id 153, "Integration IAM role")
#   off  5360646,4593939 -> 5325120,4863322
$sh = $s.Shapes.Item(7)
$sh.Left = 419.3008270515817
$sh.Top  = 382.93878175748966

# Shape 8 - "Graphic 49" picture (id 154, icon above "Integration IAM role")
#   off  5589246,4202826 -> 5553720,4453343
$sh = $s.Shapes.Item(8)
$sh.Left = 437.3008270515817
$sh.Top  = 350.65696718386516

# Shape 9 - "Rectangle 154" (id 155, "New AWS account")
#   off  4628944,3851341 -> 4671771,4136464
$sh = $s.Shapes.Item(9)
$sh.Left = 367.8560180519761
$sh.Top  = 325.70587157166364

# Shape 23 - "TextBox 16" (id 54, "Cloud One endpoint")
#   off  418688,4714893 -> 394736,4087758
$sh = $s.Shapes.Item(23)
$sh.Left = 31.08161355315006
$sh.Top  = 321.8707428114242

# Shape 24 - "Graphic 6" picture (id 55, icon above "Cloud One endpoint")
#   off  942720,4202826 -> 928965,3598238
$sh = $s.Shapes.Item(24)
$sh.Left = 73.14688874369851
$sh.Top  = 283.3258667216626

# Shape 25 - "Straight Arrow Connector 33" (id 34)
#   ext  0,914400 -> 0,1188720
$sh = $s.Shapes.Item(25)
$sh.Height = 93.60004039999774

# Shape 26 - "Straight Arrow Connector 36" (id 37)
#   off  1401800,4453771 -> 1357706,3828082
#   ext  3200400,4617 -> 4389120,4617
$sh = $s.Shapes.Item(26)
$sh.Left  = 106.90602495196624
$sh.Top   = 301.423812877566
$sh.Width = 345.60003660000825

# Shape 27 - "Straight Arrow Connector 37" (id 38)
#   off  6003891,2417787 -> 5995182,2417787
#   headEnd type="none" -> type="arrow"
$sh = $s.Shapes.Item(27)
$sh.Left = 472.06161500315864
$sh.Line.BeginArrowheadStyle = 3
# Re-touch EndArrowheadStyle (idempotent - already "arrow") so the host
# re-serialises <a:tailEnd> after <a:headEnd>, keeping the original
# headEnd/tailEnd element order instead of appending headEnd last.
$sh.Line.EndArrowheadStyle = 3
